# Generate Report for Handback
# Update the "Correspond Handoff Datetime" / "Correspond Handback DateTime" timestamps
# (and the rolled-up "Latest HO Xliff Generate Date" on the Overview sheet) for the
# 734e3018-944d-419c-a6b0-d8654137328a row after a fresh handback report generation.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G2").Value = "2016-09-06 09:20:30"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H2").Value = "2016-09-06 09:19:59"
$zhcn.Range("K2").Value = "2016-09-06 09:21:22"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H2").Value = "2016-09-06 09:20:30"
$dede.Range("K2").Value = "2016-09-06 09:21:41"
